$d = $word.ActiveDocument

$replacements = @(
  @{old="60÷9=6, 6"; new="36÷6=6, 0"},
  @{old="22÷9=2, 4"; new="56÷3=18, 2"},
  @{old="46÷4=11, 2"; new="11÷6=1, 5"},
  @{old="31÷6=5, 1"; new="64÷8=8, 0"},
  @{old="66÷9=7, 3"; new="63÷3=21, 0"},
  @{old="39÷9=4, 3"; new="64÷5=12, 4"},
  @{old="69÷4=17, 1"; new="22÷8=2, 6"},
  @{old="71÷3=23, 2"; new="68÷5=13, 3"},
  @{old="61÷9=6, 7"; new="60÷4=15, 0"},
  @{old="93÷8=11, 5"; new="88÷6=14, 4"},
  @{old="59÷3=19, 2"; new="32÷4=8, 0"},
  @{old="90÷2=45, 0"; new="97÷3=32, 1"},
  @{old="54÷4=13, 2"; new="18÷5=3, 3"},
  @{old="17÷7=2, 3"; new="72÷2=36, 0"},
  @{old="11÷2=5, 1"; new="27÷6=4, 3"},
  @{old="19÷5=3, 4"; new="45÷5=9, 0"},
  @{old="76÷9=8, 4"; new="48÷9=5, 3"},
  @{old="57÷3=19, 0"; new="37÷4=9, 1"},
  @{old="52÷4=13, 0"; new="15÷9=1, 6"},
  @{old="45÷8=5, 5"; new="51÷4=12, 3"},
  @{old="84÷2=42, 0"; new="25÷3=8, 1"},
  @{old="91÷7=13, 0"; new="15÷9=1, 6"},
  @{old="86÷5=17, 1"; new="44÷2=22, 0"},
  @{old="76÷4=19, 0"; new="37÷4=9, 1"},
  @{old="58÷5=11, 3"; new="80÷6=13, 2"}
)

foreach ($r in $replacements) {
  $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                           $true, 1, $false, $r.new, 2)
}
